$d = $word.ActiveDocument

# The document has three inline pictures living in header/footer stories:
#   - Primary (default) footer  -> Pearson logo, currently named "image2.png"
#   - First-page footer         -> Pearson logo, currently named "image2.png"
#   - First-page header         -> BTEC logo,    currently named "image1.jpg"
#
# Walk every story range and rename each inline picture based on its
# distinguishing alt-text (the "descr" attribute), which uniquely identifies
# the BTEC logo vs. the Pearson logo regardless of section/story ordering.

$storyRanges = $d.StoryRanges
foreach ($story in $storyRanges) {
    $shapes = $story.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $desc = $shp.AlternativeText

        if ($desc -eq "BTec_Logo-Orange") {
            # BTec_Logo-Orange.jpg was "image1.jpg" -> rename to "image2.jpg"
            $shp.Name = "image2.jpg"
        } elseif ($desc -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            # PearsonLogo.png was "image2.png" -> rename to "image1.png"
            $shp.Name = "image1.png"
        }
    }
}
